$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2:B145 values
$bValues = @(
1278,  1391,  1414,  1448,  1494,  1454,  1535,  1488,  1498,  1528,  1517,  1501,  1547,  1589,  1553,  1594,  1548,  1609,  1566,  1565,  1601,  1580,  1580,  1540,  1565,  1564,  1566,  1552,  1579,  1598,  1571,  1571,  1578,  1569,  1543,  1540,  1557,  1574,  1524,  1545,  1509,  1544,  1586,  1566,  1530,  1521,  1497,  1535,  1536,  1484,  1474,  1281,  1289,  1412,  1416,  1455,  1477,  1519,  1550,  1536,  1531,  1531,  1526,  1776,  1464,  1493,  1495,  1399,  1430,  1463,  1468,  1433,  1402,  1400,  1442,  1384,  1394,  1427,  1461,  1444,  1426,  1431,  1430,  1422,  1453,  1426,  1441,  1415,  1432,  1418,  1406,  1445,  1458,  1465,  1480,  1485,  1516,  1414,  1408,  1411,  1490,  1462,  1453,  1444,  1226,  1287,  1304,  1334,  1330,  1372,  1419,  1350,  1396,  1435,  1399,  1409,  1409,  1496,  1459,  1415,  1476,  1476,  1487,  1559,  1577,  1546,  1565,  1504,  1528,  1530,  1539,  1496,  1509,  1556,  1558,  1541,  1515,  1515,  1501,  1495,  1544,  1505,  1510,  1524
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $bValues[$i]
}

# Append new rows 146:180
$newRows = @(
@(44479,1496),  @(44486,1489),  @(44493,1535),  @(44500,1564),  @(44507,1495),  @(44514,1534),  @(44521,1479),  @(44528,1520),  @(44535,1521),  @(44542,1471),  @(44549,1461),  @(44556,1485),  @(44563,1197),  @(44570,1293),  @(44577,1365),  @(44584,1398),  @(44591,1460),  @(44598,1489),  @(44605,1506),  @(44612,1555),  @(44619,1546),  @(44626,1574),  @(44633,1549),  @(44640,1576),  @(44647,1594),  @(44654,1612),  @(44661,1636),  @(44668,1654),  @(44675,1610),  @(44682,1589),  @(44689,1615),  @(44696,1626),  @(44703,1653),  @(44710,1643),  @(44717,1250)
)
$startRow = 146
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $dateVal = $newRows[$i][0]
    $countVal = $newRows[$i][1]
    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $countVal
    $ws.Cells.Item($r, 3).Value = "W"
}
